$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 55557030
$ws.Range("I40").Value = 1526.3636
$ws.Range("J40").Value = 142858540
$ws.Range("K40").Value = 1526.3636
$ws.Range("L40").Value = 142858540
$ws.Range("M40").Value = -1351.3636
$ws.Range("N40").Value = -142858890

$ws.Range("H51").Value = 1975.7778
$ws.Range("I51").Value = 1900
$ws.Range("K51").Value = 1900
$ws.Range("M51").Value = -1416

$ws.Range("H86").Value = 7444.15
$ws.Range("I86").Value = 8162.8667
$ws.Range("J86").Value = 5288
$ws.Range("K86").Value = 8162.8667
$ws.Range("L86").Value = 5288
$ws.Range("M86").Value = -7039.8667
$ws.Range("N86").Value = -7534

$ws.Range("H89").Value = 7444.15
$ws.Range("I89").Value = 8162.8667
$ws.Range("J89").Value = 5288
$ws.Range("K89").Value = 40814.3335
$ws.Range("L89").Value = 26440
$ws.Range("M89").Value = -35198.3335
$ws.Range("N89").Value = -37672

$ws.Range("H137").Value = 4000
$ws.Range("I137").Value = 1875
$ws.Range("J137").Value = 4944.4443
$ws.Range("K137").Value = 5625
$ws.Range("L137").Value = 14833.3329
$ws.Range("M137").Value = -3075
$ws.Range("N137").Value = -19933.3329

$ws.Range("H138").Value = 4986.59
$ws.Range("I138").Value = 1712.0416
$ws.Range("J138").Value = 6020.6577
$ws.Range("K138").Value = 5136.1248
$ws.Range("L138").Value = 18061.9731
$ws.Range("M138").Value = 3.875200000000405
$ws.Range("N138").Value = -28341.9731

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 31948.969
$ws.Range("I32").Value = 12637.547
$ws.Range("K32").Value = 12637.547
$ws.Range("M32").Value = -12350.547

$ws.Range("H44").Value = 24669.834
$ws.Range("J44").Value = 34004.75
$ws.Range("L44").Value = 34004.75
$ws.Range("N44").Value = -34980.75

$ws.Range("H52").Value = 20496.666
$ws.Range("J52").Value = 20496.666
$ws.Range("L52").Value = 20496.666
$ws.Range("N52").Value = -21132.666

$ws.Range("H74").Value = 2028.1702
$ws.Range("I74").Value = 1276.0646
$ws.Range("J74").Value = 3485.375
$ws.Range("K74").Value = 1276.0646
$ws.Range("L74").Value = 3485.375
$ws.Range("M74").Value = -402.0645999999999
$ws.Range("N74").Value = -5233.375

$ws.Range("H77").Value = 2028.1702
$ws.Range("I77").Value = 1276.0646
$ws.Range("J77").Value = 3485.375
$ws.Range("K77").Value = 6380.322999999999
$ws.Range("L77").Value = 17426.875
$ws.Range("M77").Value = -2012.322999999999
$ws.Range("N77").Value = -26162.875

$ws.Range("H132").Value = 2815.475
$ws.Range("I132").Value = 1680.0588
$ws.Range("J132").Value = 9249.5
$ws.Range("K132").Value = 5040.1764
$ws.Range("L132").Value = 27748.5
$ws.Range("M132").Value = -2510.1764
$ws.Range("N132").Value = -32808.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H54").Value = 2229.5334
$ws.Range("I54").Value = 786.9167
$ws.Range("J54").Value = 8000
$ws.Range("K54").Value = 786.9167
$ws.Range("L54").Value = 8000
$ws.Range("M54").Value = -302.9167
$ws.Range("N54").Value = -8968

$ws.Range("H82").Value = 23166.176
$ws.Range("J82").Value = 39001.11
$ws.Range("L82").Value = 39001.11
$ws.Range("N82").Value = -39767.11

$ws.Range("H85").Value = 23166.176
$ws.Range("J85").Value = 39001.11
$ws.Range("L85").Value = 39001.11
$ws.Range("N85").Value = -41653.11

$ws.Range("H94").Value = 810.0909
$ws.Range("I94").Value = 695.1429000000001
$ws.Range("J94").Value = 1011.25
$ws.Range("K94").Value = 695.1429000000001
$ws.Range("L94").Value = 1011.25
$ws.Range("M94").Value = -244.1429000000001
$ws.Range("N94").Value = -1913.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3733.16
$ws.Range("I31").Value = 1913.9429
$ws.Range("J31").Value = 7978
$ws.Range("K31").Value = 1913.9429
$ws.Range("L31").Value = 7978
$ws.Range("M31").Value = -1618.9429
$ws.Range("N31").Value = -8568

$ws.Range("H34").Value = 3733.16
$ws.Range("I34").Value = 1913.9429
$ws.Range("J34").Value = 7978
$ws.Range("K34").Value = 1913.9429
$ws.Range("L34").Value = 7978
$ws.Range("M34").Value = -1711.9429
$ws.Range("N34").Value = -8382

$ws.Range("H50").Value = 7687
$ws.Range("J50").Value = 8954.333000000001
$ws.Range("L50").Value = 8954.333000000001
$ws.Range("N50").Value = -10204.333

$ws.Range("H51").Value = 9382.200000000001
$ws.Range("J51").Value = 9382.200000000001
$ws.Range("L51").Value = 9382.200000000001
$ws.Range("N51").Value = -10854.2

$ws.Range("H60").Value = 11090.2
$ws.Range("J60").Value = 11090.2
$ws.Range("L60").Value = 11090.2
$ws.Range("N60").Value = -12112.2

$ws.Range("H61").Value = 9382.200000000001
$ws.Range("J61").Value = 9382.200000000001
$ws.Range("L61").Value = 9382.200000000001
$ws.Range("N61").Value = -10078.2

$ws.Range("H131").Value = 49490
$ws.Range("J131").Value = 49490
$ws.Range("L131").Value = 49490
$ws.Range("N131").Value = -59570

$ws.Range("H132").Value = 2266.3096
$ws.Range("I132").Value = 1914.0968
$ws.Range("J132").Value = 3258.9092
$ws.Range("K132").Value = 5742.2904
$ws.Range("L132").Value = 9776.7276
$ws.Range("M132").Value = -3212.2904
$ws.Range("N132").Value = -14836.7276

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H10").Value = 1600.9
$ws.Range("I10").Value = 120.545456
$ws.Range("J10").Value = 3410.2222
$ws.Range("K10").Value = 361.636368
$ws.Range("L10").Value = 10230.6666
$ws.Range("M10").Value = -222.636368
$ws.Range("N10").Value = -10508.6666

$ws.Range("H11").Value = 1040.6
$ws.Range("I11").Value = 466.66666
$ws.Range("J11").Value = 1901.5
$ws.Range("K11").Value = 1399.99998
$ws.Range("L11").Value = 5704.5
$ws.Range("M11").Value = -1259.99998
$ws.Range("N11").Value = -5984.5

$ws.Range("H13").Value = 386.14285
$ws.Range("I13").Value = 375.25
$ws.Range("J13").Value = 400.66666
$ws.Range("K13").Value = 1125.75
$ws.Range("L13").Value = 1201.99998
$ws.Range("M13").Value = -957.75
$ws.Range("N13").Value = -1537.99998

$ws.Range("J17").Value = 22500.223
$ws.Range("L17").Value = 67500.66900000001
$ws.Range("N17").Value = -67838.66900000001

$ws.Range("H20").Value = 8458.777
$ws.Range("J20").Value = 8458.777
$ws.Range("L20").Value = 25376.331
$ws.Range("N20").Value = -25830.331

$ws.Range("H26").Value = 186.41667
$ws.Range("J26").Value = 105.28571
$ws.Range("L26").Value = 315.85713
$ws.Range("N26").Value = -891.85713

$ws.Range("H70").Value = 5890.1113
$ws.Range("I70").Value = 2003.4286
$ws.Range("J70").Value = 8363.454
$ws.Range("K70").Value = 6010.2858
$ws.Range("L70").Value = 25090.362
$ws.Range("M70").Value = -5695.2858
$ws.Range("N70").Value = -25720.362

$ws.Range("H73").Value = 5890.1113
$ws.Range("I73").Value = 2003.4286
$ws.Range("J73").Value = 8363.454
$ws.Range("K73").Value = 6010.2858
$ws.Range("L73").Value = 25090.362
$ws.Range("M73").Value = -4918.2858
$ws.Range("N73").Value = -27274.362

$ws.Range("H132").Value = 948956.7
$ws.Range("I132").Value = 1882199.1
$ws.Range("J132").Value = 15714.286
$ws.Range("K132").Value = 16939791.9
$ws.Range("L132").Value = 141428.574
$ws.Range("M132").Value = -16937261.9
$ws.Range("N132").Value = -146488.574

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1709.0769
$ws.Range("I102").Value = 1773.091
$ws.Range("J102").Value = 1357
$ws.Range("K102").Value = 1773.091
$ws.Range("L102").Value = 1357
$ws.Range("M102").Value = -151.0909999999999
$ws.Range("N102").Value = -4601

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2980
$ws.Range("I46").Value = 1966.6666
$ws.Range("J46").Value = 4500
$ws.Range("K46").Value = 1966.6666
$ws.Range("L46").Value = 4500
$ws.Range("M46").Value = -1778.6666
$ws.Range("N46").Value = -4876

$ws.Range("H61").Value = 3117.1333
$ws.Range("I61").Value = 2185.7
$ws.Range("K61").Value = 2185.7
$ws.Range("M61").Value = -1983.7

$ws.Range("H74").Value = 24970
$ws.Range("J74").Value = 24970
$ws.Range("L74").Value = 24970
$ws.Range("N74").Value = -26966

$ws.Range("H77").Value = 24970
$ws.Range("J77").Value = 24970
$ws.Range("L77").Value = 74910
$ws.Range("N77").Value = -84894

$ws.Range("H113").Value = 3117.1333
$ws.Range("I113").Value = 2185.7
$ws.Range("K113").Value = 2185.7
$ws.Range("M113").Value = -15.69999999999982

$ws.Range("H122").Value = 4243.2856
$ws.Range("I122").Value = 3745.818
$ws.Range("J122").Value = 4790.5
$ws.Range("K122").Value = 11237.454
$ws.Range("L122").Value = 14371.5
$ws.Range("M122").Value = -8787.454000000002
$ws.Range("N122").Value = -19271.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H109").Value = 21362
$ws.Range("J109").Value = 21362
$ws.Range("L109").Value = 21362
$ws.Range("N109").Value = -24136

$ws.Range("H122").Value = 17823.77
$ws.Range("I122").Value = 29823.428
$ws.Range("J122").Value = 3824.1667
$ws.Range("K122").Value = 89470.284
$ws.Range("L122").Value = 11472.5001
$ws.Range("M122").Value = -87020.284
$ws.Range("N122").Value = -16372.5001

$ws.Range("H126").Value = 3800.3635
$ws.Range("I126").Value = 3511.0667
$ws.Range("J126").Value = 4420.2856
$ws.Range("K126").Value = 10533.2001
$ws.Range("L126").Value = 13260.8568
$ws.Range("M126").Value = -8063.2001
$ws.Range("N126").Value = -18200.8568
